# Refresh the crypto price/volume snapshot (GitHub Actions scheduled run).
# Writes values as TEXT (matching the sheet's existing inlineStr cells) even
# when a price string happens to look numeric, by force-entering it with a
# leading apostrophe the same way a user typing into Excel would.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Text
    )
    if ($Text -match '^\s*[+-]?(\d+\.?\d*|\.\d+)\s*$') {
        $ws.Range($Cell).Value = "'" + $Text
    } else {
        $ws.Range($Cell).Value = $Text
    }
}

# Row 2 - Bitcoin
Set-TextValue "D2" "65.326.18"
Set-TextValue "E2" "  -4.83%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.253.34"
Set-TextValue "E3" "  -5.80%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.12%  "

# Row 5 - BNB
Set-TextValue "D5" "554.01"
Set-TextValue "E5" "  -3.62%  "

# Row 6 - Solana
Set-TextValue "D6" "179.29"
Set-TextValue "E6" "  -5.68%  "

# Row 7 - USDC
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.08%  "

# Row 8 - XRP
Set-TextValue "E8" "  -2.97%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "3.252.44"
Set-TextValue "E9" "  -5.51%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.183"
Set-TextValue "E10" "  -8.56%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.583"
Set-TextValue "E11" "  -4.75%  "

# Row 12 - Avalanche
Set-TextValue "D12" "47.14"
Set-TextValue "E12" "  -7.36%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000263"
Set-TextValue "E13" "  -6.96%  "

# Row 14/15 - Polkadot and BitcoinCash swapped places in the ranking
Set-TextValue "B14" "Polkadot"
Set-TextValue "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "8.51"
Set-TextValue "E14" "  -5.81%  "

Set-TextValue "B15" "BitcoinCash"
Set-TextValue "C15" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D15" "626.35"
Set-TextValue "E15" "  -1.40%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "3.771.34"
Set-TextValue "E16" "  -5.66%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "65.357.30"
Set-TextValue "E17" "  -4.52%  "

# Row 18/19 - TRON and Chainlink swapped places in the ranking
Set-TextValue "B18" "TRON"
Set-TextValue "C18" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D18" "0.116"
Set-TextValue "E18" "  -3.25%  "

Set-TextValue "B19" "Chainlink"
Set-TextValue "C19" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D19" "17.68"
Set-TextValue "E19" "  -2.23%  "

# Row 20 - WrappedEther
Set-TextValue "D20" "3.247.51"
Set-TextValue "E20" "  -6.07%  "

# Row 21 - Uniswap
Set-TextValue "D21" "11.31"
Set-TextValue "E21" "  -7.87%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.898"
Set-TextValue "E22" "  -4.08%  "

# Row 23 - InternetComputer(DFINITY)
Set-TextValue "D23" "17.69"
Set-TextValue "E23" "  -0.87%  "

# Row 24 - Litecoin
Set-TextValue "D24" "105.82"
Set-TextValue "E24" "  +6.81%  "

# Row 25 - Toncoin
Set-TextValue "D25" "4.96"
Set-TextValue "E25" "  -7.09%  "

# Row 26 - PancakeSwap
Set-TextValue "D26" "3.95"
Set-TextValue "E26" "  -7.18%  "

# Row 27 - ImmutableX
Set-TextValue "D27" "2.65"
Set-TextValue "E27" "  -5.92%  "

# Row 28 - RenderToken
Set-TextValue "D28" "9.50"
Set-TextValue "E28" "  -2.60%  "

# Row 29 - Filecoin
Set-TextValue "D29" "8.66"
Set-TextValue "E29" "  -5.59%  "

# Row 30 - EthereumClassic
Set-TextValue "D30" "30.13"
Set-TextValue "E30" "  -6.56%  "

# Row 31 - dogwifhat
Set-TextValue "D31" "4.04"
Set-TextValue "E31" "  -3.00%  "

# Row 32 - NEARProtocol
Set-TextValue "E32" "  -6.12%  "

# Row 33 - Cosmos
Set-TextValue "D33" "10.99"
Set-TextValue "E33" "  -4.67%  "

# Row 34 - Bittensor
Set-TextValue "D34" "546.11"
Set-TextValue "E34" "  +9.24%  "

# Row 35 - Hedera
Set-TextValue "E35" "  -2.78%  "

# Row 36/37 - Dai and OKB swapped places in the ranking
Set-TextValue "B36" "Dai"
Set-TextValue "C36" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D36" "0.999"
Set-TextValue "E36" "  -0.09%  "

Set-TextValue "B37" "OKB"
Set-TextValue "C37" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D37" "56.72"
Set-TextValue "E37" "  -6.97%  "

# Row 38 - Maker
Set-TextValue "D38" "3.599.81"
Set-TextValue "E38" "  -1.55%  "

# Row 39 - CoreDAO
Set-TextValue "D39" "3.67"
Set-TextValue "E39" "  +8.32%  "

# Row 40 - Stacks
Set-TextValue "D40" "3.40"
Set-TextValue "E40" "  -2.79%  "

# Row 41 - Fetch.AI
Set-TextValue "E41" "  -4.55%  "

# Row 42 - Kaspa
Set-TextValue "E42" "  -1.91%  "

# Row 43 - PEPE
Set-TextValue "D43" "0.0₃0710"
Set-TextValue "E43" "  -7.96%  "

# Row 44 - InjectiveProtocol
Set-TextValue "D44" "31.77"
Set-TextValue "E44" "  -7.12%  "

# Row 45 - TheGraph
Set-TextValue "E45" "  -7.97%  "

# Row 46 - ApeXProtocol
Set-TextValue "E46" "  -1.62%  "

# Row 47 - VeChain
Set-TextValue "D47" "0.0411"
Set-TextValue "E47" "  -5.36%  "

# Row 48 - ThetaToken
Set-TextValue "E48" "  -7.10%  "

# Row 49 - Stellar
Set-TextValue "D49" "0.128"
Set-TextValue "E49" "  -3.65%  "

# Row 50 - FirstDigitalUSD
Set-TextValue "D50" "0.999"
Set-TextValue "E50" "  -0.02%  "

# Row 51 - Mantle
Set-TextValue "E51" "  +1.82%  "
